$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "250.70")
# keep their exact textual representation instead of being coerced to a
# floating point number (which would drop trailing zeros / add fp noise).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.332.00"
$ws.Range("D3").Value = "2.094.57"
$ws.Range("D5").Value = "250.70"
$ws.Range("D8").Value = "51.72"
$ws.Range("D9").Value = "61.67"
$ws.Range("D10").Value = "0.373"
$ws.Range("D11").Value = "0.0743"
$ws.Range("D13").Value = "15.02"
$ws.Range("D14").Value = "2.403.80"
$ws.Range("D15").Value = "0.833"
$ws.Range("D16").Value = "2.103.09"
$ws.Range("D17").Value = "5.12"
$ws.Range("D18").Value = "37.287.82"
$ws.Range("D19").Value = "72.24"
$ws.Range("D20").Value = "14.12"
$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("D22").Value = "240.20"
$ws.Range("D23").Value = "5.21"
$ws.Range("D26").Value = "170.85"
$ws.Range("D27").Value = "9.22"
$ws.Range("D28").Value = "20.66"
$ws.Range("D31").Value = "1.07"
$ws.Range("D32").Value = "4.49"
$ws.Range("D34").Value = "20.73"
$ws.Range("D35").Value = "0.0917"
$ws.Range("D37").Value = "2.31"
$ws.Range("D38").Value = "1.84"
$ws.Range("D39").Value = "4.09"
$ws.Range("D41").Value = "18.45"
$ws.Range("D43").Value = "1.16"
$ws.Range("D44").Value = "99.08"
$ws.Range("D45").Value = "0.0916"
$ws.Range("D47").Value = "3.00"
$ws.Range("D48").Value = "1.319.12"
$ws.Range("D50").Value = "2.293.80"

# Restore original (default) style on column D now that the text is set,
# so only the cell contents change -- matches the source workbook which
# has no explicit style on these cells.
$dRange.Style = "Normal"

$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("E3").Value = "  +4.08%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +15.51%  "
$ws.Range("E9").Value = "  +5.58%  "
$ws.Range("E10").Value = "  +3.62%  "
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("E12").Value = "  +7.35%  "
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("E14").Value = "  +4.07%  "
$ws.Range("E15").Value = "  +4.15%  "
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("E17").Value = "  +4.91%  "
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("E20").Value = "  +9.49%  "
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("E23").Value = "  +6.55%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +5.51%  "
$ws.Range("E27").Value = "  +8.76%  "
$ws.Range("E28").Value = "  +5.14%  "
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("E31").Value = "  +28.63%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("E35").Value = "  +13.54%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +8.99%  "
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("E41").Value = "  +16.13%  "
$ws.Range("E42").Value = "  +4.59%  "
$ws.Range("E43").Value = "  +7.59%  "
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("E45").Value = "  +13.34%  "
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  +8.61%  "
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("E49").Value = "  +14.17%  "
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("E51").Value = "  +3.72%  "
